# Write intro slides, add pdf version of title, add github .png screenshot.
#
# This script applies two small text merges on the title slide:
#  1. In the "Introducing pacea: ..." textbox, the runs "of" and
#     " ecosystem information to help facilitate an ecosystem approach "
#     are merged into a single run "of ecosystem information to help
#     facilitate an ecosystem approach ".
#  2. In the "PBS / Friday 10th November 2023" textbox, the runs
#     "Friday " and "10" are merged into a single run "Friday 10".
#
# Re-assigning the .Text of a Characters() range that spans two runs
# with identical formatting causes PowerPoint to collapse them back
# into one run (keeping the first run's rPr), which matches the target
# diff exactly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. "Text Box 4": Introducing pacea: an R package of ecosystem ... ---
$titleShape = $s.Shapes.Item(2)
$titleRange = $titleShape.TextFrame.TextRange
$titleText = $titleRange.Text

$mergedPhrase = "of ecosystem information to help facilitate an ecosystem approach "
$startPos = $titleText.IndexOf("of ecosystem information to help facilitate an ecosystem approach ") + 1

$mergeChars = $titleRange.Characters($startPos, $mergedPhrase.Length)
$mergeChars.Text = $mergedPhrase

# --- 2. "Rectangle 5": PBS / Friday 10th November 2023 ---
$dateShape = $s.Shapes.Item(3)
$dateRange = $dateShape.TextFrame.TextRange
$dateText = $dateRange.Text

$mergedDate = "Friday 10"
$dateStartPos = $dateText.IndexOf("Friday 10") + 1

$dateMergeChars = $dateRange.Characters($dateStartPos, $mergedDate.Length)
$dateMergeChars.Text = $mergedDate
